# Weekly data refresh: two new price observations (fecha 44491) are added
# to the top of the data block (rows 182-183), pushing every existing
# record down by two rows (old row 182 -> new row 184, ..., old row 208 ->
# new row 210). Column headers (row 1) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 182, shifting rows 182:208
# (and their formatting) down to 184:210.
$ws.Range("A182:A183").EntireRow.Insert()

# --- New row 182: Brócoli Primera, Región Metropolitana ---------------
$ws.Cells.Item(182, 1).Value2  = 4
$ws.Cells.Item(182, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(182, 3).Value2  = "Los Lagos"
$ws.Cells.Item(182, 4).Value2  = 44491
$ws.Cells.Item(182, 5).Value2  = 10
$ws.Cells.Item(182, 6).Value2  = 100112023
$ws.Cells.Item(182, 7).Value2  = "Brócoli"
$ws.Cells.Item(182, 8).Value2  = "Sin especificar"
$ws.Cells.Item(182, 9).Value2  = "Primera"
$ws.Cells.Item(182, 10).Value2 = 700
$ws.Cells.Item(182, 11).Value2 = 1200
$ws.Cells.Item(182, 12).Value2 = 1200
$ws.Cells.Item(182, 13).Value2 = 1200
$ws.Cells.Item(182, 14).Value2 = "$/unidad"
$ws.Cells.Item(182, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(182, 16).Value2 = 1200
$ws.Cells.Item(182, 17).Value2 = 1
$ws.Cells.Item(182, 18).Value2 = "Hortaliza"

# --- New row 183: Brócoli Segunda, Región Metropolitana ----------------
$ws.Cells.Item(183, 1).Value2  = 4
$ws.Cells.Item(183, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(183, 3).Value2  = "Los Lagos"
$ws.Cells.Item(183, 4).Value2  = 44491
$ws.Cells.Item(183, 5).Value2  = 10
$ws.Cells.Item(183, 6).Value2  = 100112023
$ws.Cells.Item(183, 7).Value2  = "Brócoli"
$ws.Cells.Item(183, 8).Value2  = "Sin especificar"
$ws.Cells.Item(183, 9).Value2  = "Segunda"
$ws.Cells.Item(183, 10).Value2 = 700
$ws.Cells.Item(183, 11).Value2 = 1000
$ws.Cells.Item(183, 12).Value2 = 1000
$ws.Cells.Item(183, 13).Value2 = 1000
$ws.Cells.Item(183, 14).Value2 = "$/unidad"
$ws.Cells.Item(183, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(183, 16).Value2 = 1000
$ws.Cells.Item(183, 17).Value2 = 1
$ws.Cells.Item(183, 18).Value2 = "Hortaliza"
